# Updated cryptos list (price / volume refresh + a couple of row reorders)
# Values that look like plain numbers but must stay literal text (so Excel
# doesn't silently reparse "572.49" as a float, etc.) are written with a
# leading apostrophe, exactly like typing a quote-prefixed text value in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.440.89'
$ws.Range("E2").Value = '  -2.54%  '
$ws.Range("D3").Value = '3.179.75'
$ws.Range("E3").Value = '  -4.17%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'572.49"
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("E7").Value = '  -6.55%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '3.189.33'
$ws.Range("E9").Value = '  -3.80%  '
$ws.Range("E10").Value = '  -4.02%  '
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("D13").Value = '3.740.44'
$ws.Range("E13").Value = '  -3.92%  '
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").Value = '64.511.57'
$ws.Range("E15").Value = '  -2.51%  '
$ws.Range("D16").Value = "'25.32"
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("E17").Value = '  -3.72%  '
$ws.Range("D18").Value = '3.191.49'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("D19").Value = "'420.94"
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("E21").Value = '  -3.21%  '
$ws.Range("E22").Value = '  -2.81%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = "'70.28"
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  +2.31%  '
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("E28").Value = '  -7.69%  '
$ws.Range("D29").Value = "'8.76"
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("D30").Value = "'0.996"
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("E31").Value = '  -4.38%  '
$ws.Range("D32").Value = "'21.81"
$ws.Range("E32").Value = '  -2.49%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").Value = "'5.04"
$ws.Range("E34").Value = '  -2.61%  '
$ws.Range("E35").Value = '  -3.02%  '
$ws.Range("D36").Value = "'157.01"
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("E37").Value = '  -4.32%  '
$ws.Range("E38").Value = '  -4.77%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '2.708.66'
$ws.Range("E39").Value = '  -5.60%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = "'1.70"
$ws.Range("E40").Value = '  -4.98%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = "'24.39"
$ws.Range("E41").Value = '  -7.59%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = "'4.23"
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("D44").Value = "'0.716"
$ws.Range("E45").Value = '  -5.82%  '
$ws.Range("D46").Value = "'5.54"
$ws.Range("E46").Value = '  -6.10%  '
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("D48").Value = "'291.87"
$ws.Range("E48").Value = '  -6.66%  '
$ws.Range("E49").Value = '  -7.42%  '
$ws.Range("E50").Value = '  -5.63%  '
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = '  -0.17%  '
